$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.264813303947449
$ws.Range("B1").Value = 1.394700646400452
$ws.Range("C1").Value = 1.660661816596985
$ws.Range("D1").Value = 2.93171501159668
$ws.Range("E1").Value = -1
